$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: A1 gets the URL, B1 gets "unmatched" (per the target diff).
$ws.Range("A1").Value = "https://www.google.com"
$ws.Range("B1").Value = "unmatched"
